# Update the "Metadata" sheet (URL, Version, Date, Publisher)
$wb = $excel.ActiveWorkbook
$wsMeta = $wb.Worksheets.Item("Metadata")

$wsMeta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/parameter-definition-with-default"
$wsMeta.Range("B3").Value = "8.0.0"
$wsMeta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$wsMeta.Range("B9").Value = "LinuxForHealth Team"

# Update the "Elements" sheet: defaultValue extension Type(s) column (J5)
$wsElements = $wb.Worksheets.Item("Elements")
$wsElements.Range("J5").Value = "Extension {http://linuxforhealth.org/fhir/cdm/StructureDefinition/default-value}
"
